$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old trial rows first. This lets the shared-string table get rebuilt
# cleanly (unused old strings dropped) and new strings appended in the exact order
# they are (re)written below, matching the order cells are scanned row-by-row.
$ws.Range("A2:J4").ClearContents()

# Cells whose text is purely numeric-looking need an explicit Text format, otherwise
# Excel will silently store them as real numbers instead of the original text values.
$ws.Range("C2").NumberFormat = "@"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("J2").NumberFormat = "@"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("J3").NumberFormat = "@"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("J4").NumberFormat = "@"
$ws.Range("B5").NumberFormat = "@"

# Trial 1 (row 2)
$ws.Range("A2").Value = '(664, 549)'
$ws.Range("B2").Value = '(531, 552)'
$ws.Range("C2").Value = '133.0338302838793'
$ws.Range("D2").Value = '0:00:00.236185'
$ws.Range("E2").Value = '[[541, 554, datetime.timedelta(microseconds=77279), 2727.850365128794], [531, 552, datetime.timedelta(microseconds=160398), 26.915359635668104, -16838.956879095287]]'
$ws.Range("F2").Value = '2022-07-15 15:13:22.015893'
$ws.Range("G2").Value = '56.317654820175576'
$ws.Range("H2").Value = 'miss'
$ws.Range("I2").Value = '(645.5, 540.0)'
$ws.Range("J2").Value = '115.12710367241937'

# Trial 2 (row 3)
$ws.Range("A3").Value = '(1083, 541)'
$ws.Range("B3").Value = '(1454, 549)'
$ws.Range("C3").Value = '371.0862433451286'
$ws.Range("D3").Value = '0:00:00.500964'
$ws.Range("E3").Value = '[[1151, 549, datetime.timedelta(microseconds=77005), 637.4467985693265], [1180, 549, datetime.timedelta(microseconds=151974), 80.78136172415458, -3662.8991593639166], [1329, 559, datetime.timedelta(microseconds=234724), 269.33149240499637, 803.2844135275549], [1452, 557, datetime.timedelta(microseconds=322965), 161.24621248094928, -334.6656136858393], [1452, 552, datetime.timedelta(microseconds=413939), 5.113474851769626, -377.18779247468746]]'
$ws.Range("F3").Value = '2022-07-15 15:13:27.348214'
$ws.Range("G3").Value = '157.09317634943778'
$ws.Range("H3").Value = 'miss'
$ws.Range("I3").Value = '(1273.5, 540.0)'
$ws.Range("J3").Value = '180.7242374447877'

# Trial 3 (row 4, newly added)
$ws.Range("A4").Value = '(1038, 523)'
$ws.Range("B4").Value = '(1287, 486)'
$ws.Range("C4").Value = '251.73398658107334'
$ws.Range("D4").Value = '0:00:00.621378'
$ws.Range("E4").Value = '[[1067, 512, datetime.timedelta(microseconds=81938), 160.24505742125305], [1114, 510, datetime.timedelta(microseconds=164523), 121.04491592989531, -238.26541876429278], [1216, 507, datetime.timedelta(microseconds=250489), 172.4573631309366, 205.24832308421253], [1227, 499, datetime.timedelta(microseconds=346485), 16.618196599269826, -449.77175500141936], [1253, 494, datetime.timedelta(microseconds=442484), 25.330508243596956, 19.689551812782227], [1287, 489, datetime.timedelta(microseconds=538484), 27.016843771710054, 3.131635346849856]]'
$ws.Range("F4").Value = '2022-07-15 15:13:32.185694'
$ws.Range("G4").Value = '106.56738765265439'
$ws.Range("H4").Value = 'hit'
$ws.Range("I4").Value = '(1273.5, 540.0)'
$ws.Range("J4").Value = '55.66192594583842'

# Subject code row (row 5, newly added) -- number of trials metadata
$ws.Range("A5").Value = 'Subject Code:'
$ws.Range("B5").Value = '1'

# Restore Normal style on the numeric-text cells so no stray formatting remains
$ws.Range("C2").Style = "Normal"
$ws.Range("G2").Style = "Normal"
$ws.Range("J2").Style = "Normal"
$ws.Range("C3").Style = "Normal"
$ws.Range("G3").Style = "Normal"
$ws.Range("J3").Style = "Normal"
$ws.Range("C4").Style = "Normal"
$ws.Range("G4").Style = "Normal"
$ws.Range("J4").Style = "Normal"
$ws.Range("B5").Style = "Normal"

